# Apply cryptos.xlsx price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.955.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.11%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.086.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.67%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.086.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "

# Row 9
$ws.Range("E9").Value = "  -0.53%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.52%  "

# Row 11
$ws.Range("E11").Value = "  +0.66%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.21%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.43%  "

# Row 15
$ws.Range("E15").Value = "  +0.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.598.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.867.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18
$ws.Range("E18").Value = "  -0.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.082.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.684"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.59%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.71%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.79%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "

# Row 28
$ws.Range("E28").Value = "  -0.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.46%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.16%  "

# Row 31
$ws.Range("E31").Value = "  -2.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.25%  "

# Row 33
$ws.Range("E33").Value = "  -0.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0939"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "47.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.71%  "

# Row 37
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.945"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.16%  "

# Row 38
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.67%  "

# Row 39
$ws.Range("E39").Value = "  +2.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "48.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "

# Row 41
$ws.Range("E41").Value = "  +0.47%  "

# Row 42
$ws.Range("E42").Value = "  -0.39%  "

# Row 43
$ws.Range("E43").Value = "  +8.72%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.64%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.793.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "368.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.79%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0343"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.85%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.98%  "

# Row 49
$ws.Range("E49").Value = "  +0.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.08%  "

# Row 51
$ws.Range("E51").Value = "  +6.84%  "
